$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z1").Value = "Services"
$ws.Range("Z1").HorizontalAlignment = $ws.Range("Y1").HorizontalAlignment
$ws.Range("Z1").VerticalAlignment = $ws.Range("Y1").VerticalAlignment
$ws.Range("Z1").WrapText = $ws.Range("Y1").WrapText
$ws.Range("Z1").Font.Bold = $ws.Range("Y1").Font.Bold
$ws.Range("Z1").Font.Name = $ws.Range("Y1").Font.Name
$ws.Range("Z1").Font.Size = $ws.Range("Y1").Font.Size

